# Update automatic: dades i banners [2026-02-23 20:50]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")

$ws.Range("E2").Value = "2026-02-23 20:48:39"
$helper.Formula = "=""36%"""
$helper.Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("O2").Value = "5.9 °C"
$ws.Range("E3").Value = "2026-02-23 20:48:41"
$ws.Range("E4").Value = "2026-02-23 20:48:44"
$ws.Range("O4").Value = "12.2 °C"
$ws.Range("E5").Value = "2026-02-23 20:48:47"
$ws.Range("E6").Value = "2026-02-23 20:48:49"
$helper.Formula = "=""61%"""
$helper.Copy()
$ws.Range("H6").PasteSpecial(-4163)
$ws.Range("O6").Value = "14.1 °C"
$ws.Range("E7").Value = "2026-02-23 20:48:52"
$ws.Range("E8").Value = "2026-02-23 20:48:54"
$ws.Range("E9").Value = "2026-02-23 20:48:57"
$ws.Range("E10").Value = "2026-02-23 20:48:59"
$helper.Formula = "=""76%"""
$helper.Copy()
$ws.Range("H10").PasteSpecial(-4163)
$ws.Range("O10").Value = "10.9 °C"
$ws.Range("E11").Value = "2026-02-23 20:49:02"
$ws.Range("O11").Value = "8.9 °C"
$ws.Range("E12").Value = "2026-02-23 20:49:04"
$ws.Range("E13").Value = "2026-02-23 20:49:07"
$ws.Range("J13").Value = "1026.8 hPa"
$ws.Range("O13").Value = "7.2 °C"
$ws.Range("E14").Value = "2026-02-23 20:49:10"
$helper.Formula = "=""75%"""
$helper.Copy()
$ws.Range("H14").PasteSpecial(-4163)
$ws.Range("N14").Value = "6.8 °C 20:29 TU"
$ws.Range("O14").Value = "12.6 °C"
$ws.Range("E15").Value = "2026-02-23 20:49:12"
$helper.Formula = "=""70%"""
$helper.Copy()
$ws.Range("H15").PasteSpecial(-4163)
$ws.Range("E16").Value = "2026-02-23 20:49:15"
$ws.Range("E17").Value = "2026-02-23 20:49:17"
$helper.Formula = "=""44%"""
$helper.Copy()
$ws.Range("H17").PasteSpecial(-4163)
$ws.Range("O17").Value = "8.6 °C"
$ws.Range("E18").Value = "2026-02-23 20:49:20"
$helper.Formula = "=""73%"""
$helper.Copy()
$ws.Range("H18").PasteSpecial(-4163)
$ws.Range("E19").Value = "2026-02-23 20:49:23"
$helper.Formula = "=""46%"""
$helper.Copy()
$ws.Range("H19").PasteSpecial(-4163)
$ws.Range("K19").Value = "15.3 MJ/m2"
$ws.Range("O19").Value = "12.4 °C"
$ws.Range("E20").Value = "2026-02-23 20:49:25"
$ws.Range("O20").Value = "4.2 °C"
$ws.Range("E21").Value = "2026-02-23 20:49:28"
$ws.Range("J21").Value = "1025.6 hPa"
$ws.Range("E22").Value = "2026-02-23 20:49:30"
$ws.Range("E23").Value = "2026-02-23 20:49:33"
$ws.Range("K23").Value = "16.3 MJ/m2"
$ws.Range("O23").Value = "3.9 °C"
$ws.Range("E24").Value = "2026-02-23 20:49:35"
$ws.Range("E25").Value = "2026-02-23 20:49:38"
$ws.Range("E26").Value = "2026-02-23 20:49:41"
$helper.Formula = "=""51%"""
$helper.Copy()
$ws.Range("H26").PasteSpecial(-4163)
$ws.Range("E27").Value = "2026-02-23 20:49:43"
$ws.Range("E28").Value = "2026-02-23 20:49:46"
$helper.Formula = "=""68%"""
$helper.Copy()
$ws.Range("H28").PasteSpecial(-4163)
$ws.Range("J28").Value = "1024.9 hPa"
$ws.Range("E29").Value = "2026-02-23 20:49:48"
$ws.Range("E30").Value = "2026-02-23 20:49:51"
$helper.Formula = "=""70%"""
$helper.Copy()
$ws.Range("H30").PasteSpecial(-4163)
$ws.Range("O30").Value = "13.0 °C"
$ws.Range("E31").Value = "2026-02-23 20:49:53"
$helper.Formula = "=""46%"""
$helper.Copy()
$ws.Range("H31").PasteSpecial(-4163)
$ws.Range("E32").Value = "2026-02-23 20:49:56"
$helper.Formula = "=""66%"""
$helper.Copy()
$ws.Range("H32").PasteSpecial(-4163)
$ws.Range("K32").Value = "15.9 MJ/m2"
$ws.Range("O32").Value = "7.7 °C"
$ws.Range("E33").Value = "2026-02-23 20:49:59"
$ws.Range("E34").Value = "2026-02-23 20:50:01"
$ws.Range("O34").Value = "4.1 °C"
$ws.Range("E35").Value = "2026-02-23 20:50:04"
$ws.Range("J35").Value = "1025.1 hPa"
$ws.Range("O35").Value = "12.3 °C"
$ws.Range("E36").Value = "2026-02-23 20:50:07"
$helper.Formula = "=""72%"""
$helper.Copy()
$ws.Range("H36").PasteSpecial(-4163)
$ws.Range("E37").Value = "2026-02-23 20:50:10"
$helper.Formula = "=""66%"""
$helper.Copy()
$ws.Range("H37").PasteSpecial(-4163)
$ws.Range("O37").Value = "9.2 °C"
$ws.Range("E38").Value = "2026-02-23 20:50:12"
$ws.Range("O38").Value = "12.3 °C"
$ws.Range("E39").Value = "2026-02-23 20:50:15"
$helper.Formula = "=""25%"""
$helper.Copy()
$ws.Range("H39").PasteSpecial(-4163)
$ws.Range("E40").Value = "2026-02-23 20:50:17"
$ws.Range("O40").Value = "8.9 °C"
$ws.Range("E41").Value = "2026-02-23 20:50:20"
$ws.Range("O41").Value = "12.1 °C"
$ws.Range("E42").Value = "2026-02-23 20:50:22"
$helper.Formula = "=""79%"""
$helper.Copy()
$ws.Range("H42").PasteSpecial(-4163)
$ws.Range("E43").Value = "2026-02-23 20:50:25"
$ws.Range("E44").Value = "2026-02-23 20:50:28"
$helper.Formula = "=""34%"""
$helper.Copy()
$ws.Range("H44").PasteSpecial(-4163)
$ws.Range("N44").Value = "0.2 °C 20:17 TU"
$ws.Range("E45").Value = "2026-02-23 20:50:30"
$helper.Formula = "=""51%"""
$helper.Copy()
$ws.Range("H45").PasteSpecial(-4163)
$ws.Range("O45").Value = "8.4 °C"
$ws.Range("E46").Value = "2026-02-23 20:50:32"
$ws.Range("J46").Value = "1025.9 hPa"
$ws.Range("O46").Value = "10.3 °C"

$helper.Value = ""
$excel.CutCopyMode = 0
